$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 0: materialize & correct the built-in "Hyperlink" character style in
# styles.xml by touching it on a disposable paragraph that is deleted right
# after, so the real body of the document is left untouched by the operation.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$tempPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tempPara.Range.Text = "x"
$tempRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$tempRange.MoveEnd(1, -1)
$tempRange.Style = $d.Styles.Item("Hyperlink")

$hs = $d.Styles.Item("Hyperlink")
$hs.BaseStyle = $d.Styles.Item("Fontepargpadro")
$hs.Priority = 99
$hs.UnhideWhenUsed = 1

$d.Paragraphs.Item($d.Paragraphs.Count).Range.Delete()

# ---------------------------------------------------------------------------
# Step 1: the final (empty) paragraph of the document holds the "_GoBack"
# bookmark. Insert a new paragraph before it and fill it with the full
# sentence of body text (this lets the new run correctly inherit the
# paragraph's run-formatting, including rFonts/sz/szCs).
# ---------------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPara.Range.InsertParagraphBefore()

$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bodyRange = $bodyPara.Range
$bodyRange.Text = "Acessar o site https://fonts.google.com/ e escolher uma fonte."

# ---------------------------------------------------------------------------
# Step 2: insert another new paragraph before the body paragraph for the
# centered, bold heading "Adicionando uma fonte".
# ---------------------------------------------------------------------------
$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bodyPara.Range.InsertParagraphBefore()

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 2)
$headingRange = $headingPara.Range
$headingRange.Text = "Adicionando uma fonte"
$headingRange.ParagraphFormat.Alignment = 1
$headingRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# Step 3: turn the URL inside the body paragraph into a real hyperlink.
# ---------------------------------------------------------------------------
$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bodyRange = $bodyPara.Range
$urlText = "https://fonts.google.com/"
$urlStart = $bodyRange.Start + $bodyRange.Text.IndexOf($urlText)
$urlEnd = $urlStart + $urlText.Length
$urlRange = $d.Range($urlStart, $urlEnd)
$d.Hyperlinks.Add($urlRange, $urlText)

# Make sure every run in the paragraph carries the expected Verdana/12pt
# (sz 24 half-points) character formatting, matching the surrounding text.
$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bodyRange = $bodyPara.Range
$bodyRange.MoveEnd(1, -1)
$bodyRange.Font.Name = "Verdana"
$bodyRange.Font.Size = 12

Write-Host "Done"
